# Apply the edits described by the diff:
#  - Column header E1: "Egg Pappardelle" -> "Plain Egg Pappardelle"
#  - Cell E2: "A1PAPP" -> "01PAPP"
#  - Active selection moves from O1 to G1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Plain Egg Pappardelle"
$ws.Range("E2").Value = "01PAPP"

$ws.Range("G1").Select()
